$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '92.454.57'
Set-TextValue $ws.Range('E2') '  +6.89%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.311.73'
Set-TextValue $ws.Range('E3') '  +1.28%  '

# Row 4
Set-TextValue $ws.Range('D4') '0.998'
Set-TextValue $ws.Range('E4') '  -0.22%  '

# Row 5
Set-TextValue $ws.Range('D5') '215.67'
Set-TextValue $ws.Range('E5') '  +1.93%  '

# Row 6
Set-TextValue $ws.Range('D6') '625.70'
Set-TextValue $ws.Range('E6') '  -0.44%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.423'
Set-TextValue $ws.Range('E7') '  +13.87%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.710'
Set-TextValue $ws.Range('E8') '  +2.69%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.997'
Set-TextValue $ws.Range('E9') '  -0.24%  '

# Row 10
Set-TextValue $ws.Range('D10') '3.308.42'
Set-TextValue $ws.Range('E10') '  +1.22%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.596'
Set-TextValue $ws.Range('E11') '  +4.12%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.0000269'
Set-TextValue $ws.Range('E12') '  +5.05%  '

# Row 13
Set-TextValue $ws.Range('D13') '0.181'
Set-TextValue $ws.Range('E13') '  +1.13%  '

# Row 14
Set-TextValue $ws.Range('D14') '34.65'
Set-TextValue $ws.Range('E14') '  +1.58%  '

# Row 15
Set-TextValue $ws.Range('D15') '3.890.85'
Set-TextValue $ws.Range('E15') '  +0.46%  '

# Row 16
Set-TextValue $ws.Range('D16') '91.594.15'
Set-TextValue $ws.Range('E16') '  +5.87%  '

# Row 17
Set-TextValue $ws.Range('D17') '5.41'
Set-TextValue $ws.Range('E17') '  +1.79%  '

# Row 18
Set-TextValue $ws.Range('D18') '3.287.90'
Set-TextValue $ws.Range('E18') '  +0.48%  '

# Row 19
Set-TextValue $ws.Range('D19') '3.31'
Set-TextValue $ws.Range('E19') '  +8.64%  '

# Row 20
Set-TextValue $ws.Range('D20') '14.21'
Set-TextValue $ws.Range('E20') '  +1.48%  '

# Row 21
Set-TextValue $ws.Range('D21') '438.28'
Set-TextValue $ws.Range('E21') '  +1.64%  '

# Row 22
Set-TextValue $ws.Range('D22') '9.06'
Set-TextValue $ws.Range('E22') '  +2.18%  '

# Row 23
Set-TextValue $ws.Range('D23') '5.36'
Set-TextValue $ws.Range('E23') '  +0.70%  '

# Row 24
Set-TextValue $ws.Range('D24') '0.0000189'
Set-TextValue $ws.Range('E24') '  +45.78%  '

# Row 25
Set-TextValue $ws.Range('D25') '5.47'
Set-TextValue $ws.Range('E25') '  +7.34%  '

# Row 26
Set-TextValue $ws.Range('D26') '12.45'
Set-TextValue $ws.Range('E26') '  -0.62%  '

# Row 27
Set-TextValue $ws.Range('D27') '3.463.38'
Set-TextValue $ws.Range('E27') '  +0.49%  '

# Row 28
Set-TextValue $ws.Range('D28') '76.38'
Set-TextValue $ws.Range('E28') '  +0.40%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.00'
Set-TextValue $ws.Range('E29') '  -0.02%  '

# Row 30
Set-TextValue $ws.Range('D30') '0.181'
Set-TextValue $ws.Range('E30') '  +4.23%  '

# Row 31
Set-TextValue $ws.Range('D31') '1.00'
Set-TextValue $ws.Range('E31') '  -0.25%  '

# Row 32
Set-TextValue $ws.Range('B32') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D32') '8.78'
Set-TextValue $ws.Range('E32') '  -0.18%  '

# Row 33
Set-TextValue $ws.Range('B33') 'Bittensor'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D33') '562.47'
Set-TextValue $ws.Range('E33') '  +3.79%  '

# Row 34
Set-TextValue $ws.Range('D34') '7.45'
Set-TextValue $ws.Range('E34') '  +6.72%  '

# Row 35
Set-TextValue $ws.Range('D35') '3.70'
Set-TextValue $ws.Range('E35') '  +26.46%  '

# Row 36
Set-TextValue $ws.Range('D36') '1.35'
Set-TextValue $ws.Range('E36') '  -5.29%  '

# Row 37
Set-TextValue $ws.Range('D37') '1.93'
Set-TextValue $ws.Range('E37') '  -0.61%  '

# Row 38
Set-TextValue $ws.Range('D38') '22.75'
Set-TextValue $ws.Range('E38') '  +1.32%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.134'
Set-TextValue $ws.Range('E39') '  -2.34%  '

# Row 40
Set-TextValue $ws.Range('E40') '  +3.78%  '

# Row 41
Set-TextValue $ws.Range('E41') '  +0.03%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.397'
Set-TextValue $ws.Range('E42') '  +1.09%  '

# Row 43
Set-TextValue $ws.Range('D43') '2.00'
Set-TextValue $ws.Range('E43') '  +0.59%  '

# Row 44
Set-TextValue $ws.Range('E44') '  +0.21%  '

# Row 45
Set-TextValue $ws.Range('D45') '183.38'
Set-TextValue $ws.Range('E45') '  +2.44%  '

# Row 46
Set-TextValue $ws.Range('E46') '  -5.43%  '

# Row 47
Set-TextValue $ws.Range('B47') 'Stellar'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D47') '0.131'
Set-TextValue $ws.Range('E47') '  +6.83%  '

# Row 48
Set-TextValue $ws.Range('B48') 'OKB'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D48') '44.05'
Set-TextValue $ws.Range('E48') '  -0.74%  '

# Row 49
Set-TextValue $ws.Range('D49') '1.30'
Set-TextValue $ws.Range('E49') '  -0.14%  '

# Row 50
Set-TextValue $ws.Range('D50') '25.71'
Set-TextValue $ws.Range('E50') '  +6.90%  '

# Row 51
Set-TextValue $ws.Range('D51') '0.634'
Set-TextValue $ws.Range('E51') '  +1.23%  '
